{"js": "// Replace the date title and the 25 equation cells with their new values.\n// All old text values are unique within the document, so a direct\n// search-and-replace (matchCase, whole-text match) is safe and precise.\nconst replacements = [\n  [\"2024-01-20 Saturday\", \"2024-01-21 Sunday\"],\n  [\"20\u00d718=360\", \"66\u00d769=4554\"],\n  [\"58\u00d746=2668\", \"99\u00d736=3564\"],\n  [\"17\u00d713=221\", \"26\u00d723=598\"],\n  [\"26\u00d717=442\", \"94\u00d712=1128\"],\n  [\"64\u00d746=2944\", \"39\u00d748=1872\"],\n  [\"65\u00d786=5590\", \"53\u00d753=2809\"],\n  [\"31\u00d790=2790\", \"54\u00d782=4428\"],\n  [\"65\u00d784=5460\", \"28\u00d722=616\"],\n  [\"57\u00d760=3420\", \"46\u00d784=3864\"],\n  [\"34\u00d766=2244\", \"25\u00d787=2175\"],\n  [\"81\u00d748=3888\", \"12\u00d788=1056\"],\n  [\"50\u00d790=4500\", \"28\u00d786=2408\"],\n  [\"52\u00d711=572\", \"30\u00d727=810\"],\n  [\"48\u00d745=2160\", \"60\u00d757=3420\"],\n  [\"29\u00d789=2581\", \"32\u00d768=2176\"],\n  [\"98\u00d783=8134\", \"67\u00d779=5293\"],\n  [\"86\u00d781=6966\", \"38\u00d765=2470\"],\n  [\"11\u00d783=913\", \"75\u00d724=1800\"],\n  [\"21\u00d773=1533\", \"76\u00d739=2964\"],\n  [\"82\u00d735=2870\", \"57\u00d741=2337\"],\n  [\"55\u00d730=1650\", \"83\u00d720=1660\"],\n  [\"32\u00d760=1920\", \"14\u00d737=518\"],\n  [\"75\u00d733=2475\", \"82\u00d723=1886\"],\n  [\"88\u00d755=4840\", \"19\u00d722=418\"],\n  [\"64\u00d787=5568\", \"19\u00d764=1216\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date title and the 25 equation cells with their new values.\n# All old text values are unique within the document, so Find/Execute with\n# Replace:=wdReplaceAll (2) against the whole document's Content range is a\n# safe, precise way to apply each substitution exactly once.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-01-20 Saturday\", \"2024-01-21 Sunday\"),\n    @(\"20\u00d718=360\", \"66\u00d769=4554\"),\n    @(\"58\u00d746=2668\", \"99\u00d736=3564\"),\n    @(\"17\u00d713=221\", \"26\u00d723=598\"),\n    @(\"26\u00d717=442\", \"94\u00d712=1128\"),\n    @(\"64\u00d746=2944\", \"39\u00d748=1872\"),\n    @(\"65\u00d786=5590\", \"53\u00d753=2809\"),\n    @(\"31\u00d790=2790\", \"54\u00d782=4428\"),\n    @(\"65\u00d784=5460\", \"28\u00d722=616\"),\n    @(\"57\u00d760=3420\", \"46\u00d784=3864\"),\n    @(\"34\u00d766=2244\", \"25\u00d787=2175\"),\n    @(\"81\u00d748=3888\", \"12\u00d788=1056\"),\n    @(\"50\u00d790=4500\", \"28\u00d786=2408\"),\n    @(\"52\u00d711=572\", \"30\u00d727=810\"),\n    @(\"48\u00d745=2160\", \"60\u00d757=3420\"),\n    @(\"29\u00d789=2581\", \"32\u00d768=2176\"),\n    @(\"98\u00d783=8134\", \"67\u00d779=5293\"),\n    @(\"86\u00d781=6966\", \"38\u00d765=2470\"),\n    @(\"11\u00d783=913\", \"75\u00d724=1800\"),\n    @(\"21\u00d773=1533\", \"76\u00d739=2964\"),\n    @(\"82\u00d735=2870\", \"57\u00d741=2337\"),\n    @(\"55\u00d730=1650\", \"83\u00d720=1660\"),\n    @(\"32\u00d760=1920\", \"14\u00d737=518\"),\n    @(\"75\u00d733=2475\", \"82\u00d723=1886\"),\n    @(\"88\u00d755=4840\", \"19\u00d722=418\"),\n    @(\"64\u00d787=5568\", \"19\u00d764=1216\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n"}
